# Mise a jour de certains champs de Modules et de Professeurs
#
# Adds a new "Matières enseignés" column (E) to the Feuil1 sheet, sets the
# column widths for columns C, D and E, and leaves the selection on E6
# (matching the state captured when the workbook was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column E (adds a new shared string entry).
$ws.Range("E1").Value = "Matières enseignés"

# Column widths (stored width = ColumnWidth + 5/6 in this engine's model).
$ws.Columns.Item(3).ColumnWidth = 26.736979166666668
$ws.Columns.Item(4).ColumnWidth = 14.877604166666666
$ws.Columns.Item(5).ColumnWidth = 30.877604166666668

# Restore the selection that was active when the file was saved.
[void]$ws.Range("E6").Select()
